# Fix capitalization of "shampoo" in the "Patua & Keratin Smoothening
# Shampoo" product title. The three rows referencing this product (rows
# 11-13, column B) move from "...Smoothening Shampoo" to
# "...Smoothening shampoo" (lowercase "shampoo").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = "Patua & Keratin Smoothening shampoo"
$ws.Range("B12").Value = "Patua & Keratin Smoothening shampoo"
$ws.Range("B13").Value = "Patua & Keratin Smoothening shampoo"

# Update the saved selection/active cell (cosmetic session state).
$ws.Range("B16").Select()
